$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 7: Buchse / LUT 0820 T8CW 04 ---
$ws.Range("A7").Value = "Buchse"
$ws.Range("B7").Value = "LUT 0820 T8CW 04"
$ws.Range("C7").Value = 7.55
$ws.Range("D7").Value = 70
$ws.Range("E7").Formula = "=C7*D7"

# --- New row 8: Buchse Solar / LUM 0270-02 ---
$ws.Range("A8").Value = "Buchse Solar"
$ws.Range("B8").Value = "LUM 0270-02"
$ws.Range("C8").Value = 3.45
$ws.Range("D8").Value = 10
$ws.Range("E8").Formula = "=C8*D8"

# --- New column H: "Anzahl Stationen" header + value ---
$ws.Range("H1").Value = "Anzahl Stationen"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").NumberFormat = "General"

$ws.Range("H2").Value = 10
$ws.Range("H2").HorizontalAlignment = -4108
$ws.Range("H2").NumberFormat = "General"

$ws.Columns.Item(8).ColumnWidth = 16.25

# --- touch I1 so dimension grows to I16 like the target file ---
$ws.Range("I1").NumberFormat = "General"

# --- update the selection to match the saved cursor position ---
$ws.Range("F8").Select() | Out-Null
